$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of Fruta/hortaliza daily logic rows: reorder each row's
# date/volume/price/origin data block (D, M:T) to the new weekly snapshot.
$rows = @{
    2 = @{
        D = 44937
        M = 100
        N = 2500
        O = 3000
        P = 2750
        Q = "`$/bandeja 2 kilos"
        R = "Provincia de Diguillín"
        S = 1375
        T = 2
    }
    3 = @{
        D = 44187
        M = 80
        N = 2800
        O = 3000
        P = 2900
        Q = "`$/bandeja 2 kilos"
        R = "Provincia de Linares"
        S = 1450
        T = 2
    }
    4 = @{
        D = 44187
        M = 65
        N = 1400
        O = 1500
        P = 1446
        Q = "`$/envase 1 kilo"
        R = "Provincia de Diguillín"
        S = 1446
        T = 1
    }
    5 = @{
        D = 44942
        M = 60
        N = 2500
        O = 2500
        P = 2500
        Q = "`$/bandeja 2 kilos"
        R = "Provincia de Diguillín"
        S = 1250
        T = 2
    }
    6 = @{
        D = 44174
        M = 150
        N = 3700
        O = 3800
        P = 3747
        Q = "`$/bandeja 2 kilos"
        R = "Provincia de Linares"
        S = 1874
        T = 2
    }
    7 = @{
        D = 44596
        M = 120
        N = 2500
        O = 2700
        P = 2600
        Q = "`$/bandeja 2 kilos"
        R = "Provincia de Linares"
        S = 1300
        T = 2
    }
    8 = @{
        D = 44594
        M = 120
        N = 2500
        O = 2800
        P = 2650
        Q = "`$/bandeja 2 kilos"
        R = "Provincia de Linares"
        S = 1325
        T = 2
    }
    9 = @{
        D = 44181
        M = 65
        N = 3600
        O = 3800
        P = 3692
        Q = "`$/bandeja 2 kilos"
        R = "Provincia de Diguillín"
        S = 1846
        T = 2
    }
    10 = @{
        D = 44181
        M = 80
        N = 1800
        O = 2000
        P = 1875
        Q = "`$/envase 1 kilo"
        R = "Provincia de Diguillín"
        S = 1875
        T = 1
    }
    11 = @{
        D = 44932
        M = 60
        N = 3000
        O = 3000
        P = 3000
        Q = "`$/bandeja 2 kilos"
        R = "Provincia de Diguillín"
        S = 1500
        T = 2
    }
    12 = @{
        D = 44931
        M = 100
        N = 3000
        O = 3000
        P = 3000
        Q = "`$/bandeja 2 kilos"
        R = "Provincia de Diguillín"
        S = 1500
        T = 2
    }
    13 = @{
        D = 44935
        M = 50
        N = 3000
        O = 3000
        P = 3000
        Q = "`$/bandeja 2 kilos"
        R = "Provincia de Diguillín"
        S = 1500
        T = 2
    }
    14 = @{
        D = 44540
        M = 240
        N = 3500
        O = 3800
        P = 3650
        Q = "`$/bandeja 2 kilos"
        R = "Región del Maule"
        S = 1825
        T = 2
    }
    15 = @{
        D = 44944
        M = 60
        N = 2500
        O = 2500
        P = 2500
        Q = "`$/bandeja 2 kilos"
        R = "Provincia de Diguillín"
        S = 1250
        T = 2
    }
    16 = @{
        D = 44539
        M = 200
        N = 3800
        O = 4000
        P = 3900
        Q = "`$/bandeja 2 kilos"
        R = "Región del Maule"
        S = 1950
        T = 2
    }
}

foreach ($r in $rows.Keys) {
    $row = $rows[$r]
    $ws.Range("D$r").Value = $row.D
    $ws.Range("M$r").Value = $row.M
    $ws.Range("N$r").Value = $row.N
    $ws.Range("O$r").Value = $row.O
    $ws.Range("P$r").Value = $row.P
    $ws.Range("Q$r").Value = $row.Q
    $ws.Range("R$r").Value = $row.R
    $ws.Range("S$r").Value = $row.S
    $ws.Range("T$r").Value = $row.T
}
